# Regenerate the "K" column (column G) values for each row of save_data.
# These are the recalculated strike-count values (K) replacing the old
# Strike# values, as described in the commit message:
#   "regen save_data to use K instead of Strike#, regen std/mean,
#    calc and write s_vals"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G ("K")
$kValues = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 0
    6  = 1
    8  = 1
    9  = 1
    10 = 0
    11 = 0
    12 = 2
    13 = 1
    14 = 1
    15 = 1
    16 = 1
    18 = 1
    19 = 0
    21 = 1
    22 = 0
    23 = 2
    24 = 1
    25 = 3
    26 = 1
    28 = 3
    29 = 2
    30 = 2
    31 = 2
    32 = 3
    33 = 2
    34 = 1
    35 = 1
    36 = 3
    37 = 2
    38 = 0
    39 = 1
    40 = 1
    43 = 1
    44 = 4
    45 = 1
    46 = 2
    47 = 1
    48 = 0
    49 = 1
    50 = 1
    51 = 0
    52 = 1
    53 = 0
    54 = 2
    55 = 2
    56 = 0
    57 = 6
    58 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
